$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, B value, D text)
$updates = @{
    "Total Hits" = @(
        @{ Row = 2; B = 1421; D = "36.13%" },
        @{ Row = 3; B = 2836; D = "36.05%" },
        @{ Row = 4; B = 4215; D = "35.72%" },
        @{ Row = 5; B = 5614; D = "35.69%" },
        @{ Row = 6; B = 7044; D = "35.82%" }
    )
    "Hits_entity" = @(
        @{ Row = 2; B = 831;  D = "34.17%" },
        @{ Row = 3; B = 1648; D = "33.88%" },
        @{ Row = 4; B = 2478; D = "33.96%" },
        @{ Row = 5; B = 3317; D = "34.10%" },
        @{ Row = 6; B = 4147; D = "34.10%" }
    )
    "Hits_numerical" = @(
        @{ Row = 2; B = 169; D = "25.80%" },
        @{ Row = 3; B = 334; D = "25.50%" },
        @{ Row = 4; B = 485; D = "24.68%" },
        @{ Row = 5; B = 629; D = "24.01%" },
        @{ Row = 6; B = 809; D = "24.70%" }
    )
    "Hits_boolean" = @(
        @{ Row = 2; B = 308;  D = "53.75%" },
        @{ Row = 3; B = 631;  D = "55.06%" },
        @{ Row = 4; B = 924;  D = "53.75%" },
        @{ Row = 5; B = 1228; D = "53.58%" },
        @{ Row = 6; B = 1533; D = "53.51%" }
    )
    "Hits_date" = @(
        @{ Row = 2; B = 109; D = "41.13%" },
        @{ Row = 3; B = 217; D = "40.94%" },
        @{ Row = 4; B = 320; D = "40.25%" },
        @{ Row = 5; B = 428; D = "40.38%" },
        @{ Row = 6; B = 539; D = "40.68%" }
    )
    "Hits_string" = @(
        @{ Row = 3; B = 6;  D = "37.50%" },
        @{ Row = 4; B = 8;  D = "33.33%" },
        @{ Row = 5; B = 12; D = "37.50%" },
        @{ Row = 6; B = 16; D = "40.00%" }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 2).Value = $entry.B
        $dCell = $ws.Cells.Item($entry.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $entry.D
    }
}
